$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells
$ws.Range("I1").Value = "Course Name"
$ws.Range("J1").Value = "Student ID"

# Fill new column data for rows 2-7
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 9).Value = "Associate Degree of Information Technology"
    $ws.Cells.Item($r, 10).Value = "A00123456"
}

# Match the header cell style (s="1") used by neighboring cells
$ws.Range("A1").Copy()
$ws.Range("I1:J7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update selection to match the recorded view state
$ws.Range("I2").Select()
